# Update cryptos list (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "62.710.67"

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.441.48"
$ws.Range("E3").Value = "  +1.63%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.08%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'566.82"
$ws.Range("E5").Value = "  +0.81%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'145.83"
$ws.Range("E6").Value = "  +2.51%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.09%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  +0.37%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  +2.62%  "

# Row 10 - TRON
$ws.Range("E10").Value = "  +0.57%  "

# Row 11 - Toncoin
$ws.Range("D11").Value = "'5.32"
$ws.Range("E11").Value = "  +1.14%  "

# Row 12 - Cardano
$ws.Range("D12").Value = "'0.354"
$ws.Range("E12").Value = "  +1.32%  "

# Row 13 - Avalanche
$ws.Range("D13").Value = "'26.99"
$ws.Range("E13").Value = "  +5.87%  "

# Row 14 - ShibaInu
$ws.Range("E14").Value = "  +6.18%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "2.882.99"
$ws.Range("E15").Value = "  +1.59%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "62.443.63"
$ws.Range("E16").Value = "  +0.99%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.426.52"
$ws.Range("E17").Value = "  +1.59%  "

# Row 18 - Chainlink
$ws.Range("D18").Value = "'11.28"
$ws.Range("E18").Value = "  +0.63%  "

# Row 19 - Uniswap
$ws.Range("E19").Value = "  +1.87%  "

# Row 20 - BitcoinCash
$ws.Range("D20").Value = "'323.75"
$ws.Range("E20").Value = "  +0.85%  "

# Row 21 - Polkadot
$ws.Range("D21").Value = "'4.18"
$ws.Range("E21").Value = "  +1.10%  "

# Row 23 - now SuiNetwork (was Litecoin)
$ws.Range("B23").Value = "SuiNetwork"
$ws.Range("C23").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D23").Value = "'1.83"
$ws.Range("E23").Value = "  +5.08%  "

# Row 24 - now Litecoin (was SuiNetwork)
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "'67.32"
$ws.Range("E24").Value = "  +1.83%  "

# Row 25 - Aptos
$ws.Range("D25").Value = "'8.72"
$ws.Range("E25").Value = "  -0.65%  "

# Row 26 - PEPE
$ws.Range("E26").Value = "  +9.32%  "

# Row 27 - Bittensor
$ws.Range("D27").Value = "'577.23"
$ws.Range("E27").Value = "  +2.61%  "

# Row 28 - WrappedeETH
$ws.Range("D28").Value = "2.561.65"
$ws.Range("E28").Value = "  +1.64%  "

# Row 29 - now InternetComputer(DFINITY) (was Binance-PegBSC-USD)
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").Value = "'8.44"
$ws.Range("E29").Value = "  +3.47%  "

# Row 30 - now Binance-PegBSC-USD (was InternetComputer(DFINITY))
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  -0.19%  "

# Row 31 - Fetch.AI
$ws.Range("E31").Value = "  +4.14%  "

# Row 32 - Kaspa
$ws.Range("E32").Value = "  +0.70%  "

# Row 33 - PancakeSwap
$ws.Range("E33").Value = "  +0.50%  "

# Row 34 - ImmutableX
$ws.Range("D34").Value = "'1.55"
$ws.Range("E34").Value = "  +3.36%  "

# Row 35 - NEARProtocol
$ws.Range("E35").Value = "  +3.97%  "

# Row 36 - FirstDigitalUSD
$ws.Range("E36").Value = "  -0.15%  "

# Row 37 - PolygonEcosystemToken
$ws.Range("E37").Value = "  +1.42%  "

# Row 38 - RenderToken
$ws.Range("D38").Value = "'5.42"
$ws.Range("E38").Value = "  +0.22%  "

# Row 39 - EthereumClassic
$ws.Range("E39").Value = "  +1.48%  "

# Row 40 - Monero
$ws.Range("D40").Value = "'148.21"
$ws.Range("E40").Value = "  -2.02%  "

# Row 41 - Stacks
$ws.Range("E41").Value = "  +1.98%  "

# Row 42 - USDe
$ws.Range("E42").Value = "  +0.43%  "

# Row 43 - dogwifhat
$ws.Range("E43").Value = "  +7.56%  "

# Row 44 - Aave
$ws.Range("D44").Value = "'148.45"
$ws.Range("E44").Value = "  +0.76%  "

# Row 45 - Filecoin
$ws.Range("E45").Value = "  +1.87%  "

# Row 46 - Hedera
$ws.Range("D46").Value = "'0.0535"
$ws.Range("E46").Value = "  +1.18%  "

# Row 47 - InjectiveProtocol
$ws.Range("D47").Value = "'20.53"
$ws.Range("E47").Value = "  +3.73%  "

# Row 48 - Mantle
$ws.Range("E48").Value = "  +2.71%  "

# Row 49 - VeChain
$ws.Range("E49").Value = "  +3.11%  "

# Row 50 - Stellar
$ws.Range("E50").Value = "  +0.99%  "

# Row 51 - BitgetToken
$ws.Range("E51").Value = "  +3.47%  "
